# Daily attendance processing - reorder "Recorded By" names so entries
# that start with "System, " become "<other>, System" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column G holds "Recorded By" values (header in row 1).
$col = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value = $rest + ", System"
    }
}
